# "Update Rizka - 15 Juni 2020"
#
# User activity reconstructed from the saved workbook state:
#   1. Worked on "User - Ganti Password" sheet: deleted the (now stale)
#      row 14 - the rows below shift up, and left the selection at E21.
#   2. Switched to the "User - Forgot Password" sheet and left the
#      selection at J20 - this sheet is the active tab when the file
#      was saved.

$wb = $excel.ActiveWorkbook

# --- Sheet "User - Ganti Password" (4th tab) ---
$wsGanti = $wb.Worksheets.Item(4)
$null = $wsGanti.Select()
$null = $wsGanti.Rows(14).Delete()
$null = $wsGanti.Range("E21").Select()

# --- Sheet "User - Forgot Password" (2nd tab) becomes the active sheet ---
$wsForgot = $wb.Worksheets.Item(2)
$null = $wsForgot.Activate()
$null = $wsForgot.Range("J20").Select()
